# Update the public EPEX Spot prices workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column ("21-nov") right before
#     the existing "01-oct." column (column DV), shifting every later
#     column one to the right (DV..EZ -> DW..FA).
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns("DV").Insert()
$wsPrix.Range("DV1").Value = "21-nov"
$wsPrix.Range("DV2:DV25").Value = "-"

# --- Sheet "Gaz": append the new daily row.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A155").NumberFormat = "@"
$wsGaz.Range("A155").Value = "2025-11-19"
$wsGaz.Range("A155").Style = "Normal"
$wsGaz.Range("B155").Value = 29.925

# --- Sheet "CO2": append the new daily row.
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A155").NumberFormat = "@"
$wsCo2.Range("A155").Value = "2025-11-19"
$wsCo2.Range("A155").Style = "Normal"
$wsCo2.Range("B155").Value = 80.34
